# NDA template edit: fill in COPLASIMON / Jonathan Leloux placeholders,
# drop now-unneeded highlighting, tidy the signature block.

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = $oldText
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Replacement.Text = $newText
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# 1) "BY AND BETWEEN" party-identification paragraph.
Replace-Text `
    "…………………hereinafter referred to as “XXXX“, a legal entity organised and existing under the laws of   , with its registered office at           , represented by Mr. / Mrs. ………………………...;" `
    "COPLASIMON hereinafter referred to as “COPLASIMON“, represented by Mr. Jonathan Leloux;"

# 2) "XXXX is …." whereas clause.
Replace-Text "XXXX is …." "COPLASIMON is a collaborative platform created in the framework of SerendiPV, a four years project which has received funding from the European Union's Horizon 2020 research and innovation programme under grant agreement No 953016."

# 3) "XXXX posses-ses certain information ... field of ................." whereas clause.
Replace-Text `
    "XXXX possesses certain information, including technical or business information, in relation to the field of …………………………." `
    "COPLASIMON possesses certain information, including technical or business information, in relation to the field of photovoltaics"

# Remove leftover lightGray highlight on the "one (1) year" / "five (5) years" terms.
Replace-Text "one (1) year" "one (1) year"
Replace-Text "five (5) years" "five (5) years"

# Strip the lightGray highlight left on those two runs (Find/Replace keeps formatting by
# default, so clear explicitly via the resulting range).
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "one (1) year"
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
while ($rng.Find.Execute()) {
    $rng.HighlightColorIndex = 0
    $rng.Collapse(0)
}

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "five (5) years"
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
while ($rng.Find.Execute()) {
    $rng.HighlightColorIndex = 0
    $rng.Collapse(0)
}

# Governing law paragraph: remove the yellow highlight on every run (and on pPr mark).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*governed by and construed*Brussels*") {
        $target = $p
    }
}
if ($target -ne $null) {
    $target.Range.HighlightColorIndex = 0
}

Write-Host "done-part-1"
